$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..21) down to (3..22), working from the bottom up
# so we don't overwrite rows before they've been copied.
for ($r = 21; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("A$dst").Value2 = $ws.Range("A$src").Value2
    $ws.Range("B$dst").Value2 = $ws.Range("B$src").Value2
    $ws.Range("C$dst").Value2 = $ws.Range("C$src").Value2
}

# New row inserted at row 2
$ws.Range("A2").Value2 = 0.0610952377319335
$ws.Range("B2").Value2 = 0.5822855234146118
$ws.Range("C2").Value2 = 0.008422106504440301

# New rows appended at the end (23..31)
$newRows = @(
    @(-3.810809135437012, 1.403007388114929, 0.0495486259460449),
    @(-1.585423946380615, 2.060841083526612, -2.507726192474365),
    @(-5.486822128295898, 2.457437515258789, -1.076503276824951),
    @(3.813155174255371, -5.157403945922852, 7.194998264312744),
    @(-3.507768154144287, 2.501498937606812, 0.7795240879058838),
    @(0.2215757369995117, -0.4009582996368408, 2.163901329040528),
    @(0.1625576019287109, 1.34720504283905, -0.6319388151168823),
    @(0.044438362121582, -0.1398162841796875, -0.8414495587348938),
    @(-0.1983919143676757, -0.413076639175415, 0.2017757892608642)
)

$row = 23
foreach ($vals in $newRows) {
    $ws.Range("A$row").Value2 = $vals[0]
    $ws.Range("B$row").Value2 = $vals[1]
    $ws.Range("C$row").Value2 = $vals[2]
    $row++
}
